# Insert a new row at position 56 (pushes existing rows 56-103 down to 57-104)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new record's data
$ws.Cells.Item(56, 1).Value  = 10
$ws.Cells.Item(56, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(56, 3).Value  = "La Araucanía"
$ws.Cells.Item(56, 4).Value  = 44957
$ws.Cells.Item(56, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 5).Value  = 9
$ws.Cells.Item(56, 6).Value  = 100112030
$ws.Cells.Item(56, 7).Value  = "Poroto granado"
$ws.Cells.Item(56, 8).Value  = "Sin especificar"
$ws.Cells.Item(56, 9).Value  = "Primera"
$ws.Cells.Item(56, 10).Value = 30
$ws.Cells.Item(56, 11).Value = 45000
$ws.Cells.Item(56, 12).Value = 45000
$ws.Cells.Item(56, 13).Value = 45000
$ws.Cells.Item(56, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(56, 15).Value = "Región del Maule"
$ws.Cells.Item(56, 16).Value = 1800
$ws.Cells.Item(56, 17).Value = 25
$ws.Cells.Item(56, 18).Value = "Hortaliza"
